$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, shifting rows 16:75 down to 17:76
$ws.Rows.Item(16).Insert()

# Update B15 value (row 15 stays in place)
$ws.Range("B15").Value = 4733.5

# Fill the newly inserted row 16 with its data
$ws.Range("A16").Value = 21
$ws.Range("B16").Value = 119.13
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2025
$ws.Range("E16").Value = "06/2025"
